$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (stored as serial date number)
$ws.Range("A1").Value = 45436

# Update prices in column D
$ws.Range("D29").Value = 348.194
$ws.Range("D30").Value = 368.347
$ws.Range("D31").Value = 396.64
